$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 5
    6  = 2
    7  = 8
    8  = 3
    9  = 3
    10 = 3
    11 = 7
    12 = 1
    13 = 7
    14 = 3
    15 = 6
    16 = 2
    17 = 9
    18 = 6
    19 = 7
    20 = 8
    21 = 3
    22 = 2
    23 = 2
    24 = 8
    25 = 6
    26 = 1
    27 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
